$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 32.333332
$ws.Range("I11").Value = 32.333332
$ws.Range("K11").Value = 32.333332
$ws.Range("M11").Value = 107.666668
$ws.Range("H17").Value = 1190
$ws.Range("J17").Value = 1190
$ws.Range("L17").Value = 3570
$ws.Range("N17").Value = -3906
$ws.Range("H64").Value = 3861
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3861
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3861
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4357
$ws.Range("H67").Value = 3861
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3861
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3861
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -5577
$ws.Range("H112").Value = 2254.4
$ws.Range("I112").Value = 916.3333
$ws.Range("K112").Value = 2748.9999
$ws.Range("M112").Value = -1640.9999
$ws.Range("H116").Value = 4128.2144
$ws.Range("I116").Value = 3298.3333
$ws.Range("K116").Value = 3298.3333
$ws.Range("M116").Value = 143.6667000000002
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H138").Value = 3998.625
$ws.Range("I138").Value = 2147.25
$ws.Range("K138").Value = 6441.75
$ws.Range("M138").Value = -1301.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 485.66666
$ws.Range("I4").Value = 485.66666
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 485.66666
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -369.66666
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H32").Value = 4281408
$ws.Range("I32").Value = 4379038
$ws.Range("J32").Value = 3500367.5
$ws.Range("K32").Value = 4379038
$ws.Range("L32").Value = 3500367.5
$ws.Range("M32").Value = -4378751
$ws.Range("N32").Value = -3500941.5
$ws.Range("H110").Value = 3701123.5
$ws.Range("I110").Value = 4112247
$ws.Range("J110").Value = 1013
$ws.Range("K110").Value = 4112247
$ws.Range("L110").Value = 1013
$ws.Range("M110").Value = -4110202
$ws.Range("N110").Value = -5103
$ws.Range("H122").Value = 29999.818
$ws.Range("I122").Value = 32599.8
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 97799.39999999999
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -95349.39999999999
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4219.8096
$ws.Range("I134").Value = 4214.533
$ws.Range("K134").Value = 12643.599
$ws.Range("M134").Value = -10108.599

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H10").Value = 450
$ws.Range("I10").Value = 450
$ws.Range("K10").Value = 450
$ws.Range("M10").Value = -311
$ws.Range("H31").Value = 1125.5714
$ws.Range("I31").Value = 961.3333
$ws.Range("K31").Value = 961.3333
$ws.Range("M31").Value = -666.3333
$ws.Range("H34").Value = 1125.5714
$ws.Range("I34").Value = 961.3333
$ws.Range("K34").Value = 961.3333
$ws.Range("M34").Value = -759.3333
$ws.Range("H39").Value = 27832
$ws.Range("I39").Value = 10000
$ws.Range("J39").Value = 32290
$ws.Range("K39").Value = 10000
$ws.Range("L39").Value = 32290
$ws.Range("M39").Value = -9609
$ws.Range("N39").Value = -33072
$ws.Range("H49").Value = 27832
$ws.Range("I49").Value = 10000
$ws.Range("J49").Value = 32290
$ws.Range("K49").Value = 10000
$ws.Range("L49").Value = 32290
$ws.Range("M49").Value = -9818
$ws.Range("N49").Value = -32654
$ws.Range("H134").Value = 2232.5
$ws.Range("I134").Value = 2232.5
$ws.Range("K134").Value = 6697.5
$ws.Range("M134").Value = -4162.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H7").Value = 1061.125
$ws.Range("I7").Value = 1289.8182
$ws.Range("J7").Value = 558
$ws.Range("K7").Value = 3869.4546
$ws.Range("L7").Value = 1674
$ws.Range("M7").Value = -3757.4546
$ws.Range("N7").Value = -1898
$ws.Range("H14").Value = 194.75
$ws.Range("I14").Value = 194.75
$ws.Range("K14").Value = 584.25
$ws.Range("M14").Value = -411.25
$ws.Range("H68").Value = 2973.6667
$ws.Range("J68").Value = 3224.7585
$ws.Range("L68").Value = 9674.2755
$ws.Range("N68").Value = -11296.2755
$ws.Range("H71").Value = 2973.6667
$ws.Range("J71").Value = 3224.7585
$ws.Range("L71").Value = 29022.8265
$ws.Range("N71").Value = -37134.8265
$ws.Range("H109").Value = 799
$ws.Range("I109").Value = 799
$ws.Range("K109").Value = 2397
$ws.Range("M109").Value = -1357
$ws.Range("H113").Value = 1310.8334
$ws.Range("J113").Value = 1310.8334
$ws.Range("L113").Value = 3932.5002
$ws.Range("N113").Value = -8272.5002
$ws.Range("H131").Value = 478181.16
$ws.Range("I131").Value = 1384.625
$ws.Range("K131").Value = 4153.875
$ws.Range("M131").Value = 886.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H80").Value = 4098
$ws.Range("J80").Value = 4000
$ws.Range("L80").Value = 4000
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 4098
$ws.Range("J83").Value = 4000
$ws.Range("L83").Value = 20000
$ws.Range("N83").Value = -29984
$ws.Range("H102").Value = 2067.2
$ws.Range("J102").Value = 1450
$ws.Range("L102").Value = 1450
$ws.Range("N102").Value = -4694

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1170.6666
$ws.Range("I22").Value = 737.75
$ws.Range("J22").Value = 1387.125
$ws.Range("K22").Value = 737.75
$ws.Range("L22").Value = 1387.125
$ws.Range("M22").Value = -442.75
$ws.Range("N22").Value = -1977.125
$ws.Range("H27").Value = 1170.6666
$ws.Range("I27").Value = 737.75
$ws.Range("J27").Value = 1387.125
$ws.Range("K27").Value = 737.75
$ws.Range("L27").Value = 1387.125
$ws.Range("M27").Value = -630.75
$ws.Range("N27").Value = -1601.125
$ws.Range("H46").Value = 1693.1818
$ws.Range("I46").Value = 1693.1818
$ws.Range("K46").Value = 1693.1818
$ws.Range("M46").Value = -1505.1818
$ws.Range("H61").Value = 1716.8334
$ws.Range("I61").Value = 1575.75
$ws.Range("K61").Value = 1575.75
$ws.Range("M61").Value = -1373.75
$ws.Range("H113").Value = 1716.8334
$ws.Range("I113").Value = 1575.75
$ws.Range("K113").Value = 1575.75
$ws.Range("M113").Value = 594.25
$ws.Range("H122").Value = 6225.8823
$ws.Range("I122").Value = 4528.5
$ws.Range("J122").Value = 7414.05
$ws.Range("K122").Value = 13585.5
$ws.Range("L122").Value = 22242.15
$ws.Range("M122").Value = -11135.5
$ws.Range("N122").Value = -27142.15
$ws.Range("H136").Value = 1769.375
$ws.Range("I136").Value = 1692.5
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5077.5
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2527.5
$ws.Range("N136").Value = -11100
$ws.Range("H140").Value = 59999.5
$ws.Range("J140").Value = 59999.5
$ws.Range("L140").Value = 59999.5
$ws.Range("N140").Value = -70359.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1349
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1349
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1349
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1807
$ws.Range("H32").Value = 20000
$ws.Range("I32").Value = 20000
$ws.Range("K32").Value = 20000
$ws.Range("M32").Value = -19683
$ws.Range("H62").Value = 13566.667
$ws.Range("I62").Value = 14950
$ws.Range("K62").Value = 14950
$ws.Range("M62").Value = -14326
$ws.Range("H65").Value = 13566.667
$ws.Range("I65").Value = 14950
$ws.Range("K65").Value = 74750
$ws.Range("M65").Value = -71630
$ws.Range("H94").Value = 10880
$ws.Range("J94").Value = 10880
$ws.Range("L94").Value = 10880
$ws.Range("N94").Value = -12682
$ws.Range("H113").Value = 340
$ws.Range("I113").Value = 298.4
$ws.Range("J113").Value = 392
$ws.Range("K113").Value = 895.1999999999999
$ws.Range("L113").Value = 1176
$ws.Range("M113").Value = 1274.8
$ws.Range("N113").Value = -5516
$ws.Range("H122").Value = 3639.3635
$ws.Range("I122").Value = 1894.8
$ws.Range("J122").Value = 5093.1665
$ws.Range("K122").Value = 5684.4
$ws.Range("L122").Value = 15279.4995
$ws.Range("M122").Value = -3234.4
$ws.Range("N122").Value = -20179.4995
$ws.Range("H132").Value = 2129.5
$ws.Range("I132").Value = 2161.875
$ws.Range("K132").Value = 6485.625
$ws.Range("M132").Value = -3955.625
$ws.Range("H136").Value = 2730.3684
$ws.Range("I136").Value = 2781.7778
$ws.Range("J136").Value = 1805
$ws.Range("K136").Value = 8345.3334
$ws.Range("L136").Value = 5415
$ws.Range("M136").Value = -5795.3334
$ws.Range("N136").Value = -10515

